$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.060.98'
$ws.Range('E2').Value = '  -1.51%  '
$ws.Range('D3').Value = '2.303.86'
$ws.Range('E3').Value = '  -1.88%  '
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').Value = "'316.07"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.91%  '
$ws.Range('D6').Value = "'104.68"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.54%  '
$ws.Range('D7').Value = "'0.628"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.40%  '
$ws.Range('E8').Value = '  -0.09%  '
$ws.Range('D9').Value = "'0.609"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.44%  '
$ws.Range('D10').Value = "'39.75"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -4.03%  '
$ws.Range('D11').Value = "'0.0907"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.89%  '
$ws.Range('D12').Value = "'8.44"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.22%  '
$ws.Range('E13').Value = '  +1.19%  '
$ws.Range('D14').Value = "'0.976"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.63%  '
$ws.Range('D15').Value = "'15.47"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.80%  '
$ws.Range('D16').Value = '2.651.92'
$ws.Range('D17').Value = '2.298.83'
$ws.Range('E17').Value = '  -2.11%  '
$ws.Range('D18').Value = '42.016.28'
$ws.Range('E18').Value = '  -1.51%  '
$ws.Range('E19').Value = '  -0.09%  '
$ws.Range('E20').Value = '  -0.48%  '
$ws.Range('D21').Value = "'286.29"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +11.05%  '
$ws.Range('D22').Value = "'73.84"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -4.18%  '
$ws.Range('E23').Value = '  -2.12%  '
$ws.Range('E24').Value = '  -0.48%  '
$ws.Range('D25').Value = "'10.01"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +6.79%  '
$ws.Range('D26').Value = "'1.01"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.65%  '
$ws.Range('D27').Value = "'3.98"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.00%  '
$ws.Range('D28').Value = "'10.95"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -3.67%  '
$ws.Range('D29').Value = "'23.52"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.37%  '
$ws.Range('E30').Value = '  +0.37%  '
$ws.Range('D31').Value = "'165.47"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -5.32%  '
$ws.Range('D32').Value = "'35.50"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.03%  '
$ws.Range('D33').Value = "'0.0884"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E34').Value = '  -1.02%  '
$ws.Range('D35').Value = "'5.90"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.41%  '
$ws.Range('E36').Value = '  +1.40%  '
$ws.Range('E37').Value = '  -4.84%  '
$ws.Range('D38').Value = "'4.66"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.37%  '
$ws.Range('D39').Value = "'2.93"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +8.65%  '
$ws.Range('E40').Value = '  -2.41%  '
$ws.Range('D41').Value = "'3.63"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.65%  '
$ws.Range('D42').Value = "'102.88"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +20.92%  '
$ws.Range('E43').Value = '  +2.07%  '
$ws.Range('D44').Value = "'70.66"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.72%  '
$ws.Range('E45').Value = '  -3.84%  '
$ws.Range('E46').Value = '  -0.04%  '
$ws.Range('D47').Value = "'117.11"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.28%  '
$ws.Range('D48').Value = "'12.08"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.35%  '
$ws.Range('B49').Value = 'ordi'
$ws.Range('C49').Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range('D49').Value = "'78.17"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +6.48%  '
$ws.Range('B50').Value = 'FraxShare'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D50').Value = "'9.12"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.08%  '
